$d = $word.ActiveDocument

# 1) "tool-chain decisions" -> "toolchain decisions" (spelling fix)
$d.Content.Find.Execute("tool-chain decisions worth `$100k", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "toolchain decisions worth `$100k", 2)

# 2) Merge the "<space>" + "simplifying " + "data structures. " runs into a single
#    run reading " simplifying data structures. " (same visible text, fewer runs).
$d.Content.Find.Execute(" simplifying data structures. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " simplifying data structures. ", 2)

# 3) Flip w:overflowPunct from false to true on the Normal / TOC Heading / No Spacing
#    paragraph styles (w:styleId Normal / ContentsHeading / NoSpacing).
foreach ($styleName in @("Normal", "TOC Heading", "No Spacing")) {
    $style = $d.Styles($styleName)
    $style.ParagraphFormat.HangingPunctuation = $true
}
